$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Withdrawals" to "Sheet1"
$ws.Name = "Sheet1"

# Update the "withdrawal" -> "withdraw" values in column C (rows 2-5)
$ws.Range("C2").Value = "withdraw"
$ws.Range("C3").Value = "withdraw"
$ws.Range("C4").Value = "withdraw"
$ws.Range("C5").Value = "withdraw"
